# Update visitor/attendance counts (column F) on three sheets to match
# the re-generated site data (gh-pages output at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1163
$ws1.Range("F3").Value = 589
$ws1.Range("F6").Value = 162
$ws1.Range("F7").Value = 60
$ws1.Range("F8").Value = 62
$ws1.Range("F10").Value = 5370
$ws1.Range("F11").Value = 4843

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 77

# --- Sheet "全部类型" (All types, combined) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1163
$ws4.Range("F3").Value = 589
$ws4.Range("F6").Value = 162
$ws4.Range("F7").Value = 60
$ws4.Range("F8").Value = 62
$ws4.Range("F10").Value = 5370
$ws4.Range("F11").Value = 4843
$ws4.Range("F17").Value = 77
